# Add an "environ" flag column of 1s for each data row (rows 2-7),
# mirroring the rest of the sheet's per-row data layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E7").Value = 1

# Move the active selection to E8, just below the newly filled column.
$ws.Range("E8").Select()
